$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column order P..AF matching the existing "Accuracy sweep" block.
$cols = @("P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF")

# Update the existing row 6 accuracy figure (AF6).
$ws.Range("AF6").Value = 0.22623995566638899

# New rows 7-10, appended below the existing data block.
$rowsData = @{
    7  = @(100,2,12,1,1,0,0,1,100,$true, 0.1,6,200,5,7,3000000,0.22268403066407999)
    8  = @(100,2,12,1,1,0,0,1,100,$false,0.1,6,200,5,8,3000000,0.22268403066407999)
    9  = @(100,2,12,1,1,0,0,1,100,$true, 0.1,6,200,5,7,3000000,0.22268403066407999)
    10 = @(100,2,12,1,1,0,0,1,100,$true, 0.1,6,200,5,7,3000000,0.22268403066407999)
}

foreach ($r in 7..10) {
    $vals = $rowsData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}

# Move the active selection to Y9, matching the committed view state.
$ws.Range("Y9").Select()

Write-Output "done"
